$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z10").Value = 43568
$ws.Range("Z10").NumberFormat = $ws.Range("B5").NumberFormat
